# Daily attendance processing - 2025-11-10 10:51:38
#
# The "Recorded By" column (G) lists the people/accounts that touched each
# attendance record, separated by ", ". This pass normalizes the ordering of
# that list so that the first two entries are kept in ascending ordinal
# (character-code) order, e.g.:
#   "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#   "system, System, backup@backdoor.com"     -> "System, system, backup@backdoor.com"
#   "dnasr281@gmail.com, admin@admin.com"     -> "admin@admin.com, dnasr281@gmail.com"
# Any additional trailing entries are left exactly where they are, and rows
# whose first two entries are already in order are left untouched.

function Test-OrdinalGreater($a, $b) {
    $lenA = $a.Length
    $lenB = $b.Length
    $minLen = $lenA
    if ($lenB -lt $minLen) { $minLen = $lenB }
    for ($i = 0; $i -lt $minLen; $i++) {
        $codeA = [int][char]$a[$i]
        $codeB = [int][char]$b[$i]
        if ($codeA -gt $codeB) { return $true }
        if ($codeA -lt $codeB) { return $false }
    }
    return ($lenA -gt $lenB)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) { continue }

    $parts = $text -split ', '
    if ($parts.Length -ge 2) {
        if (Test-OrdinalGreater $parts[0] $parts[1]) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $cell.Value = [string]::Join(', ', $parts)
        }
    }
}
